$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 20, 0.1,  0.25, 10, 9.5, 1, 50, 30, 60, 1, 20, 0.05),
    @(1, 20, 0.1,  0.25, 8,  9.5, 1, 50, 30, 60, 1, 20, 0.05),
    @(1, 20, 0.12, 0.25, 10, 9.5, 1, 50, 30, 60, 1, 20, 0.05),
    @(1, 20, 0.12, 0.28, 10, 9.5, 1, 50, 30, 60, 1, 20, 0.05),
    @(1, 30, 0.12, 0.28, 10, 9.5, 1, 50, 30, 60, 1, 20, 0.05)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
